# Auto-generated edit script updating cryptos list values
# Applies the latest snapshot of price (D) and 1h volume change (E)
# figures, plus two coin-name/link swaps (rows 31/32 and 37/38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '66.983.39'
$ws.Cells.Item(2, 5).Value = '  +0.03%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.494.15'
$ws.Cells.Item(3, 5).Value = '  +0.51%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.17%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '606.39'
$ws.Cells.Item(5, 5).Value = '  +0.93%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '145.62'
$ws.Cells.Item(6, 5).Value = '  -1.37%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '3.494.85'
$ws.Cells.Item(7, 5).Value = '  +0.58%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  -0.23%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.477'
$ws.Cells.Item(9, 5).Value = '  -1.27%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.141'
$ws.Cells.Item(10, 5).Value = '  -0.52%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '7.96'
$ws.Cells.Item(11, 5).Value = '  +6.81%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.418'
$ws.Cells.Item(12, 5).Value = '  -1.53%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000213'
$ws.Cells.Item(13, 5).Value = '  +0.31%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '4.077.47'
$ws.Cells.Item(14, 5).Value = '  +0.26%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '31.18'
$ws.Cells.Item(15, 5).Value = '  -1.31%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.488.53'
$ws.Cells.Item(16, 5).Value = '  +0.23%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '66.981.57'
$ws.Cells.Item(17, 5).Value = '  -0.09%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.117'
$ws.Cells.Item(18, 5).Value = '  +0.04%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '10.73'
$ws.Cells.Item(19, 5).Value = '  +8.01%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.30'
$ws.Cells.Item(20, 5).Value = '  -2.12%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '15.37'
$ws.Cells.Item(21, 5).Value = '  +0.71%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '429.36'
$ws.Cells.Item(22, 5).Value = '  -2.22%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.604'
$ws.Cells.Item(23, 5).Value = '  -2.35%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '79.31'
$ws.Cells.Item(24, 5).Value = '  +0.71%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.03%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.620.33'
$ws.Cells.Item(26, 5).Value = '  -0.09%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.0000117'
$ws.Cells.Item(27, 5).Value = '  -1.18%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.73'
$ws.Cells.Item(28, 5).Value = '  -1.23%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '8.14'
$ws.Cells.Item(29, 5).Value = '  -2.75%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.51'
$ws.Cells.Item(30, 5).Value = '  +1.16%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Fetch.AI'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.56'
$ws.Cells.Item(31, 5).Value = '  -3.11%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.00'
$ws.Cells.Item(32, 5).Value = '  -0.24%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.166'
$ws.Cells.Item(33, 5).Value = '  +0.87%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '25.37'
$ws.Cells.Item(34, 5).Value = '  -0.09%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.78'
$ws.Cells.Item(35, 5).Value = '  -1.68%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.00%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'NEARProtocol'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '5.73'
$ws.Cells.Item(37, 5).Value = '  -6.05%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Aptos'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '7.95'
$ws.Cells.Item(38, 5).Value = '  +0.79%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.998'
$ws.Cells.Item(39, 5).Value = '  -0.16%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '175.15'
$ws.Cells.Item(40, 5).Value = '  +0.87%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0894'
$ws.Cells.Item(41, 5).Value = '  +0.81%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.33'
$ws.Cells.Item(42, 5).Value = '  -0.75%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.893'
$ws.Cells.Item(43, 5).Value = '  +0.14%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.98'
$ws.Cells.Item(44, 5).Value = '  -11.11%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '46.24'
$ws.Cells.Item(45, 5).Value = '  -0.51%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '27.92'
$ws.Cells.Item(46, 5).Value = '  -6.56%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.22'
$ws.Cells.Item(47, 5).Value = '  -2.71%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '7.34'
$ws.Cells.Item(48, 5).Value = '  -2.18%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.40'
$ws.Cells.Item(49, 5).Value = '  -1.07%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.978'
$ws.Cells.Item(50, 5).Value = '  -0.85%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.245'
$ws.Cells.Item(51, 5).Value = '  -0.21%  '
